# Steps 1-3 is reviwed
#
# The "Subtitle 2" placeholder shape (an empty subtitle left over on the
# title-style slides) is removed from the two slides that still had it.

$p = $ppt.ActivePresentation

# Slide 1 ("Step 1" title slide) - drop the empty Subtitle placeholder (id 3).
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item("Subtitle 2").Delete()

# Slide 3 ("Step 2" title slide) - drop the empty Subtitle placeholder (id 3).
$slide3 = $p.Slides.Item(3)
$slide3.Shapes.Item("Subtitle 2").Delete()
